$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source country list for the June 2014 export report was re-run (new
# values + two countries inserted: the last entry "ZONA FRANCA PACIFICO -
# CLO" drops out of the ranking, so row 174 goes away and the table shrinks
# to A1:C173).
$ws.Rows(174).Delete()

# Columns: row, country (A), Sector Agro. Expo FOB US$ (B), Participacion % (C)
$data = @(
    ,@(5, "ESTADOS UNIDOS", 2266948412.38, 37.42)
    ,@(6, "VENEZUELA", 371029804.35, 6.12)
    ,@(7, "BELGICA", 351126626.37, 5.8)
    ,@(8, "REINO UNIDO", 313469123.39, 5.17)
    ,@(9, "ALEMANIA", 263524865.34, 4.35)
    ,@(10, "JAPON", 256666851.56, 4.24)
    ,@(11, "PAISES BAJOS - HOLANDA", 239168010.48, 3.95)
    ,@(12, "CANADA", 199502704.25, 3.29)
    ,@(13, "ITALIA", 161858521.38, 2.67)
    ,@(14, "PERU", 156707697.95, 2.59)
    ,@(15, "ECUADOR", 141421277.68, 2.33)
    ,@(16, "ESPANA", 128908280.91, 2.13)
    ,@(17, "CHILE", 125872890.95, 2.08)
    ,@(18, "MEXICO", 94031108.07, 1.55)
    ,@(19, "RUSIA", 87085461.6, 1.44)
    ,@(20, "COREA (SUR) REPUBLICA DE", 68570953.44, 1.13)
    ,@(21, "ZONA FRANCA CARTAGENA", 62571349.52, 1.03)
    ,@(22, "FRANCIA", 50784152.99, 0.84)
    ,@(23, "FINLANDIA", 44438949.57, 0.73)
    ,@(24, "BRASIL", 38066294.07, 0.63)
    ,@(25, "SUECIA", 36941641.33, 0.61)
    ,@(26, "PANAMA", 32853113.08, 0.54)
    ,@(27, "HAITI", 31448204.38, 0.52)
    ,@(28, "AUSTRALIA", 31031136.46, 0.51)
    ,@(29, "ZONA FRANCA BIOCOMBUSTIBLE DL CARIBE S A", 28801295.59, 0.48)
    ,@(30, "NORUEGA", 28594498, 0.47)
    ,@(31, "REPUBLICA DOMINICANA", 27732622.78, 0.46)
    ,@(32, "PUERTO RICO", 24168494.23, 0.4)
    ,@(33, "LIBANO", 19877828.66, 0.33)
    ,@(34, "TRINIDAD Y TOBAGO", 17656480.56, 0.29)
    ,@(35, "ISRAEL", 15419267.25, 0.25)
    ,@(36, "COSTA RICA", 14932891.29, 0.25)
    ,@(37, "MALASYA", 14906921.55, 0.25)
    ,@(38, "CHINA", 14125942.2, 0.23)
    ,@(39, "EMIRATOS ARABES", 13968026.1, 0.23)
    ,@(40, "JAMAICA", 13390016.8, 0.22)
    ,@(41, "ANGOLA", 11346649.69, 0.19)
    ,@(42, "BOLIVIA", 10330502.77, 0.17)
    ,@(43, "COSTA DE MARFIL", 10280591, 0.17)
    ,@(44, "ANTILLAS HOLANDESAS", 10208904.19, 0.17)
    ,@(45, "DINAMARCA", 9860162.42, 0.16)
    ,@(46, "POLONIA", 9368610.67, 0.15)
    ,@(47, "HONG KONG", 9356787.91, 0.15)
    ,@(48, "ARGENTINA", 8430952.16, 0.14)
    ,@(49, "ARGELIA", 8298644.26, 0.14)
    ,@(50, "MAURITANIA", 8063110.75, 0.13)
    ,@(51, "NUEVA ZELANDIA", 8060311.05, 0.13)
    ,@(52, "TAIWAN (FORMOSA)", 7667577.79, 0.13)
    ,@(53, "SURINAM", 7379963.13, 0.12)
    ,@(54, "RUMANIA", 6976899.07, 0.12)
    ,@(55, "GUATEMALA", 6706158.23, 0.11)
    ,@(56, "ZAIRE", 6491550, 0.11)
    ,@(57, "SRI LANKA", 6354507.67, 0.1)
    ,@(58, "TURQUIA", 5481604.53, 0.09)
    ,@(59, "SUDAFRICA REPUBLICA DE", 5168214.92, 0.09)
    ,@(60, "GRECIA", 5023659.67, 0.08)
    ,@(61, "UCRANIA", 4674359.14, 0.08)
    ,@(62, "PORTUGAL", 4549229.19, 0.08)
    ,@(63, "ARUBA", 4544958.13, 0.08)
    ,@(64, "ZONA FRANCA BARRANQUILLA", 4512961.61, 0.07)
    ,@(65, "HONDURAS", 4254384.58, 0.07)
    ,@(66, "ESLOVENIA", 4062117.71, 0.07)
    ,@(67, "REPUBLICA CHECA", 3426530.47, 0.06)
    ,@(68, "SUIZA", 3196772.42, 0.05)
    ,@(69, "BULGARIA", 3193115.64, 0.05)
    ,@(70, "NIGERIA", 3171050.23, 0.05)
    ,@(71, "MARRUECOS", 2935519.08, 0.05)
    ,@(72, "GUYANA", 2837869.49, 0.05)
    ,@(73, "EL SALVADOR", 2812769.05, 0.05)
    ,@(74, "ESTONIA", 2782846.33, 0.05)
    ,@(75, "GUADALUPE", 2711381.29, 0.04)
    ,@(76, "LIBIA (INCLUYE FEZZAN)", 2708690.97, 0.04)
    ,@(77, "BENIN", 2573257.64, 0.04)
    ,@(78, "CUBA", 2459421.04, 0.04)
    ,@(79, "ZONA FRANCA BOGOTA", 2244928.93, 0.04)
    ,@(80, "VIETNAM", 1928400.75, 0.03)
    ,@(81, "JORDANIA", 1833715.48, 0.03)
    ,@(82, "ARABIA SAUDITA", 1709291.49, 0.03)
    ,@(83, "TANZANIA REPUBLICA UNIDA DE", 1665847.49, 0.03)
    ,@(84, "SINGAPUR", 1573180.19, 0.03)
    ,@(85, "KUWAIT", 1562185.43, 0.03)
    ,@(86, "URUGUAY", 1495705.77, 0.02)
    ,@(87, "IRLANDA (EIRE)", 1486124.03, 0.02)
    ,@(88, "INDIA", 1445701.99, 0.02)
    ,@(89, "MARTINICA", 1438682.47, 0.02)
    ,@(90, "MOZAMBIQUE", 1438000.5, 0.02)
    ,@(91, "BAHAMAS", 1423547.76, 0.02)
    ,@(92, "REP GEORGIA=GEORGIA", 1211615.09, 0.02)
    ,@(93, "BELARUS - BIELORUSIA", 1105746.73, 0.02)
    ,@(94, "BARBADOS", 1069628.7, 0.02)
    ,@(95, "SIRIA REPUBLICA ARABE DE", 986139.9, 0.02)
    ,@(96, "SANTA LUCIA", 938054.86, 0.02)
    ,@(97, "THAILANDIA", 892439.54, 0.01)
    ,@(98, "BELICE", 884691.89, 0.01)
    ,@(99, "BAHRAIN", 857907.91, 0.01)
    ,@(100, "SENEGAL", 848163.2, 0.01)
    ,@(101, "ZAMBIA", 807884, 0.01)
    ,@(102, "EGIPTO", 787067.79, 0.01)
    ,@(103, "LETONIA", 733227.5, 0.01)
    ,@(104, "GABON", 698176.3, 0.01)
    ,@(105, "KENIA", 697255.64, 0.01)
    ,@(106, "ESLOVAQUIA", 662805.86, 0.01)
    ,@(107, "CONGO", 659611.33, 0.01)
    ,@(108, "HUNGRIA", 619678.47, 0.01)
    ,@(109, "PARAGUAY", 577091.17, 0.01)
    ,@(110, "SAN VICENTE Y LAS GRANADINAS", 568994.98, 0.01)
    ,@(111, "YEMEN", 566239, 0.01)
    ,@(112, "TOGO", 552380.02, 0.01)
    ,@(113, "GUINEA", 548581.45, 0.01)
    ,@(114, "DOMINICA", 507264.69, 0.01)
    ,@(115, "CAMERUN", 491169.48, 0.01)
    ,@(116, "ANTIGUA Y BARBUDA", 477744.4, 0.01)
    ,@(117, "MONACO", 423000, 0.01)
    ,@(118, "LITUANIA", 406628.61, 0.01)
    ,@(119, "AUSTRIA", 376391.37, 0.01)
    ,@(120, "INDONESIA", 366737.91, 0.01)
    ,@(121, "SAN CRISTOBAL Y NIEVES", 362900.28, 0.01)
    ,@(122, "GRANADA", 354990.66, 0.01)
    ,@(123, "LIBERIA", 354723.13, 0.01)
    ,@(124, "CROACIA", 340163.79, 0.01)
    ,@(125, "UGANDA", 326793.6, 0.01)
    ,@(126, "TUNICIA - TUNEZ", 325854.8, 0.01)
    ,@(127, "KAZAJSTAN", 320130.36, 0.01)
    ,@(128, "CABO VERDE", 308348.75, 0.01)
    ,@(129, "SIERRA LEONA", 283974.86, $null)
    ,@(130, "SUDAN", 251680, $null)
    ,@(131, "GUAYANA FRANCESA", 247938.52, $null)
    ,@(132, "FILIPINAS", 239302.11, $null)
    ,@(133, "KIRGUIZISTAN", 208476.1, $null)
    ,@(134, "ZF PERMANENTE SANTANDER - FLORIDABLANCA", 191347.2, $null)
    ,@(135, "IRAQ", 184491.1, $null)
    ,@(136, "GHANA", 183989.4, $null)
    ,@(137, "MAURICIO", 176573.29, $null)
    ,@(138, "FIJI - FIDJI", 164768.81, $null)
    ,@(139, "ISLANDIA", 160581.36, $null)
    ,@(140, "SERBIA Y MONTENEGRO", 159807.06, $null)
    ,@(141, "COMORAS", 125611.9, $null)
    ,@(142, "ZONA FRANCA LA CANDELARIA - CTG", 119371.26, $null)
    ,@(143, "UZBEKISTAN", 102661.18, $null)
    ,@(144, "ALBANIA", 101510.75, $null)
    ,@(145, "MACEDONIA", 95225.8, $null)
    ,@(146, "PAKISTAN", 88921.4, $null)
    ,@(147, "MALAWI", 87261, $null)
    ,@(148, "QATAR", 84854.46, $null)
    ,@(149, "AZERBAIJAN", 72482.35, $null)
    ,@(150, "MACAO", 68922.18, $null)
    ,@(151, "VIRGENES ISLAS (BRITANICAS)", 62524.03, $null)
    ,@(152, "ARMENIA", 55986.33, $null)
    ,@(153, "GUAM", 52088.39, $null)
    ,@(154, "NIGER", 50387, $null)
    ,@(155, "T. A. DE PALESTINA", 45200, $null)
    ,@(156, "GUINEA ECUATORIAL", 44138.6, $null)
    ,@(157, "NICARAGUA", 36023.67, $null)
    ,@(158, "ZONA FRANCA PERMANENTE LA CAYENA", 27788.39, $null)
    ,@(159, "ZONA FRANCA PERMANENTE PARQ INDL DEXTON", 25933.79, $null)
    ,@(160, "OMAN", 21542.42, $null)
    ,@(161, "ISLAS CAIMAN", 20671.5, $null)
    ,@(162, "NUEVA CALEDONIA", 18714.75, $null)
    ,@(163, "MONSERRAT ISLA", 16457.6, $null)
    ,@(164, "MARSHALL ISLAS", 8163.9, $null)
    ,@(165, "SEYCHELLES", 6860, $null)
    ,@(166, "ZONA FRANCA RIONEGRO - MEDELLIN", 3659, $null)
    ,@(167, "MOLDAVIA", 3519.23, $null)
    ,@(168, "CHIPRE", 3197.86, $null)
    ,@(169, "PALAU ISLAS", 2944.6, $null)
    ,@(170, "MALDIVAS", 2875.02, $null)
    ,@(171, "TURKMENISTAN", 2042.93, $null)
    ,@(172, "TADJIKISTAN", 1046.12, $null)
    ,@(173, "ZONA FRANCA PALMASECA - CALI", 53, $null)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    if ($item[3] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $item[3]
    } else {
        $ws.Cells.Item($r, 3).ClearContents()
    }
}
